# Bai 18 Chia cum van ban - cap nhat ten bai tap
# Title placeholder on slide 43 ("Bai tap") -> split into two runs:
#   "Bai " + "tap 18.1"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(43)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$tr.Text = "Bài "
[void]$tr.InsertAfter("tập 18.1")
